# Update the verification matrix on Sheet1 to match the reorganized
# content (rows re-sorted/re-grouped, a few helper rows removed, and
# new rows added for recently introduced functions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 'change_model/getKcatAcrossIsoenzymes.m'
$ws.Range("B9").ClearContents()
$ws.Range("A10").Value = 'gather_kcats/fuzzyKcatMatching.m'
$ws.Range("B10").Value = 'tc0010'
$ws.Range("A11").Value = 'gather_kcats/getStandardKcat.m'
$ws.Range("A12").Value = 'gather_kcats/mergeDlkcatAndFuzzyKcats.m'
$ws.Range("B12").Value = 'tc0011'
$ws.Range("A13").Value = 'gather_kcats/readDLKcatOutput.m'
$ws.Range("B13").ClearContents()
$ws.Range("A14").Value = 'gather_kcats/runDLKcat.m'
$ws.Range("B14").ClearContents()
$ws.Range("A15").Value = 'gather_kcats/selectKcatValue.m'
$ws.Range("B15").Value = 'tc0011 - we did not test all possible parameterizations here'
$ws.Range("A16").Value = 'gather_kcats/writeDLKcatInput.m'
$ws.Range("B16").Value = 'partly by tc0011'
$ws.Range("A17").Value = 'get_enzyme_data/calculateMW.m'
$ws.Range("B17").Value = 'Not explicitly tested - The code has existed for a long time'
$ws.Range("A18").Value = 'get_enzyme_data/findECInDB.m'
$ws.Range("B18").Value = 'Called from getECfromDatabase - tested when that function is tested.'
$ws.Range("A19").Value = 'get_enzyme_data/getECfromDatabase.m'
$ws.Range("B19").Value = 'tc0007 - does not test download of the databases - this is tested in the manual workflows for Yeast-GEM and Human-GEM'
$ws.Range("A20").Value = 'get_enzyme_data/getECfromGEM.m'
$ws.Range("B20").Value = 'tc0006'
$ws.Range("A21").Value = 'get_enzyme_data/getECstring.m'
$ws.Range("B21").Value = 'Called from findECInDB - tested when that function is tested.'
$ws.Range("A22").Value = 'get_enzyme_data/loadBRENDAdata.m'
$ws.Range("B22").Value = 'Called from fuzzyKcatMatching - tested when that function is tested.'
$ws.Range("A23").Value = 'get_enzyme_data/loadDatabases.m'
$ws.Range("B23").Value = 'Called from getECfromDatabase - tested when that function is tested.'
$ws.Range("A24").Value = 'kcat_sensitivity_analysis/Bayesian/*'
$ws.Range("B24").ClearContents()
$ws.Range("A25").Value = 'kcat_sensitivity_analysis/findMaxValue.m'
$ws.Range("C25").ClearContents()
$ws.Range("A26").Value = 'kcat_sensitivity_analysis/findTopLimitations.m'
$ws.Range("B26").Value = 'Not explicitly tested - The code has existed for a long time'
$ws.Range("A27").Value = 'kcat_sensitivity_analysis/sensitivityTuning.m'
$ws.Range("B27").ClearContents()
$ws.Range("A28").Value = 'kcat_sensitivity_analysis/sigmaFitter.m'
$ws.Range("C28").ClearContents()
$ws.Range("A29").Value = 'kcat_sensitivity_analysis/topUsedEnzymes.m'
$ws.Range("C29").Value = 'Not sure this works anymore - written for the old model structure?'
$ws.Range("A30").Value = 'kcat_sensitivity_analysis/truncateValues.m'
$ws.Range("B30").Value = 'Not explicitly tested - The code has existed for a long time'
$ws.Range("A31").Value = 'limit_proteins/constrainProtConcs.m'
$ws.Range("A32").Value = 'limit_proteins/findLimitingUBs.m'
$ws.Range("A33").Value = 'limit_proteins/fitGAM.m'
$ws.Range("A34").Value = 'limit_proteins/flexibilizeProteins.m'
$ws.Range("A35").Value = 'limit_proteins/measureAbundance.m'
$ws.Range("A36").Value = 'limit_proteins/readProteomics.m'
$ws.Range("A37").Value = 'limit_proteins/updateProtPool.m'
$ws.Range("A38").Value = 'model_adapter/ModelAdapter.m'
$ws.Range("B38").Value = 'Base class for adapters, not explicitly tested.'
$ws.Range("A39").Value = 'model_adapter/ModelAdapterManager.m'
$ws.Range("B39").Value = 'tc0008'
$ws.Range("A40").Value = 'utilities/*'
$ws.Range("B40").Value = 'Not explicitly tested - The code has existed for a long time. There are a lot of functions in here, unclear which of them are still useful and functional.'
$ws.Range("A41").Value = 'utilities/loadConventionalGEM.m'
$ws.Range("B41").ClearContents()
$ws.Range("A42").Value = 'utilities/loadEcModel.m'
$ws.Range("B42").ClearContents()
$ws.Range("A43").Value = 'utilities/saveEcModel.m'
$ws.Range("A44").Value = 'userdata/ecHumanGEM/HumanGEMAdapter.m'
$ws.Range("A45").Value = 'userdata/ecYeastGEM/YeastGEMAdapter.m'
$ws.Range("B45").ClearContents()
$ws.Range("A46").Value = 'userdata/ecYeastGEM/code/changeMedia_batch.m'
$ws.Range("B46").Value = 'Not explicitly tested - The code has existed for a long time'
$ws.Range("A47").Value = 'userdata/ecYeastGEM/code/getModelParameters_obsolete.m'
$ws.Range("B47").Value = 'Not used'
$ws.Range("A48").Value = 'userdata/ecYeastGEM/code/manualModifications.m'
$ws.Range("A49").Value = 'userdata/ecYeastGEM/code/removeIncorrectPathways.m'
$ws.Range("A50").Value = 'userdata/ecYeastGEM/code/scaleBioMass.m'
$ws.Range("A51").Value = 'userdata/ecYeastGEM/code/sumBioMass.m'
$ws.Range("A52").Value = 'userdata/ecYeastGEM/code/sumProtein.m'
$ws.Range("B52").Value = 'Not explicitly tested - The code has existed for a long time'

# Match the author's final cursor position recorded in the sheet view.
$ws.Range("A12").Select()
